$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# A8: "Volume 30   Number  30" -> "Volume 30   Number  31" (issue number increments)
$ws.Range("A8").Value = "Volume 30   Number  31"

# C9: report date range shifts forward by one week
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Crime-statistics table numeric updates ---
$ws.Range("F14").Value = 5
$ws.Range("H14").Value = -50
$ws.Range("I14").Value = 36
$ws.Range("J14").Value = 31
$ws.Range("K14").Value = 16.129032258064
$ws.Range("L14").Value = -10
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = -82.439024390243
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -71.428571428571
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 22
$ws.Range("H15").Value = -63.636363636363
$ws.Range("I15").Value = 78
$ws.Range("J15").Value = 118
$ws.Range("K15").Value = -33.898305084745
$ws.Range("L15").Value = -9.302325581395
$ws.Range("M15").Value = -37.096774193548
$ws.Range("N15").Value = -73.825503355704
$ws.Range("C16").Value = 39
$ws.Range("D16").Value = 55
$ws.Range("E16").Value = -29.090909090909
$ws.Range("F16").Value = 162
$ws.Range("G16").Value = 196
$ws.Range("H16").Value = -17.346938775510
$ws.Range("I16").Value = 1124
$ws.Range("J16").Value = 1246
$ws.Range("K16").Value = -9.791332263242
$ws.Range("L16").Value = 12.175648702594
$ws.Range("M16").Value = -17.352941176470
$ws.Range("N16").Value = -80.179862458120
$ws.Range("C17").Value = 53
$ws.Range("D17").Value = 61
$ws.Range("E17").Value = -13.114754098360
$ws.Range("F17").Value = 260
$ws.Range("G17").Value = 274
$ws.Range("H17").Value = -5.109489051094
$ws.Range("I17").Value = 1795
$ws.Range("J17").Value = 1867
$ws.Range("K17").Value = -3.856454204606
$ws.Range("L17").Value = 9.853121175030
$ws.Range("M17").Value = 54.608096468561
$ws.Range("N17").Value = -50
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = -21.875
$ws.Range("F18").Value = 109
$ws.Range("G18").Value = 132
$ws.Range("H18").Value = -17.424242424242
$ws.Range("I18").Value = 912
$ws.Range("J18").Value = 1081
$ws.Range("K18").Value = -15.633672525439
$ws.Range("L18").Value = 13.151364764268
$ws.Range("M18").Value = 14
$ws.Range("N18").Value = -86.007977907333
$ws.Range("C19").Value = 158
$ws.Range("D19").Value = 121
$ws.Range("E19").Value = 30.578512396694
$ws.Range("F19").Value = 538
$ws.Range("G19").Value = 552
$ws.Range("H19").Value = -2.536231884057
$ws.Range("I19").Value = 3690
$ws.Range("J19").Value = 3846
$ws.Range("K19").Value = -4.056162246489
$ws.Range("L19").Value = 29.428270782181
$ws.Range("M19").Value = 35.761589403973
$ws.Range("N19").Value = -43.195812807881
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = -36.363636363636
$ws.Range("F20").Value = 102
$ws.Range("G20").Value = 126
$ws.Range("H20").Value = -19.047619047619
$ws.Range("I20").Value = 814
$ws.Range("J20").Value = 789
$ws.Range("K20").Value = 3.168567807351
$ws.Range("L20").Value = 51.301115241635
$ws.Range("M20").Value = 145.180722891566
$ws.Range("N20").Value = -85.466880914122
$ws.Range("C21").Value = 299
$ws.Range("D21").Value = 310
$ws.Range("E21").Value = -3.548387096774
$ws.Range("F21").Value = 1184
$ws.Range("G21").Value = 1312
$ws.Range("H21").Value = -9.756097560975
$ws.Range("I21").Value = 8449
$ws.Range("J21").Value = 8978
$ws.Range("K21").Value = -5.892180886611
$ws.Range("L21").Value = 21.446025585741
$ws.Range("M21").Value = 29.486590038314
$ws.Range("N21").Value = -70.227985482222
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 66.666666666666
$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 26
$ws.Range("H22").Value = -46.153846153846
$ws.Range("I22").Value = 168
$ws.Range("J22").Value = 190
$ws.Range("K22").Value = -11.578947368421
$ws.Range("L22").Value = 27.272727272727
$ws.Range("M22").Value = 21.739130434782
$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 20
$ws.Range("E23").Value = 45
$ws.Range("F23").Value = 109
$ws.Range("G23").Value = 111
$ws.Range("H23").Value = -1.801801801801
$ws.Range("I23").Value = 747
$ws.Range("J23").Value = 773
$ws.Range("K23").Value = -3.363518758085
$ws.Range("L23").Value = 3.319502074688
$ws.Range("M23").Value = 52.760736196319
$ws.Range("C24").Value = 294
$ws.Range("D24").Value = 283
$ws.Range("E24").Value = 3.886925795053
$ws.Range("F24").Value = 1136
$ws.Range("G24").Value = 1222
$ws.Range("H24").Value = -7.037643207855
$ws.Range("I24").Value = 8461
$ws.Range("J24").Value = 9399
$ws.Range("K24").Value = -9.979785083519
$ws.Range("L24").Value = 18.951216083227
$ws.Range("M24").Value = 56.772280896794
$ws.Range("C25").Value = 89
$ws.Range("D25").Value = 84
$ws.Range("E25").Value = 5.952380952380
$ws.Range("F25").Value = 363
$ws.Range("G25").Value = 308
$ws.Range("H25").Value = 17.857142857142
$ws.Range("I25").Value = 2747
$ws.Range("J25").Value = 2740
$ws.Range("K25").Value = 0.255474452554
$ws.Range("L25").Value = 14.220374220374
$ws.Range("M25").Value = -16.096518020769
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -53.125
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 184
$ws.Range("K26").Value = -20.108695652173
$ws.Range("L26").Value = -1.342281879194
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = -53.846153846153
$ws.Range("F27").Value = 49
$ws.Range("G27").Value = 51
$ws.Range("H27").Value = -3.921568627450
$ws.Range("I27").Value = 367
$ws.Range("J27").Value = 414
$ws.Range("K27").Value = -11.352657004830
$ws.Range("L27").Value = -3.166226912928
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("I28").Value = 100
$ws.Range("J28").Value = 121
$ws.Range("K28").Value = -17.355371900826
$ws.Range("L28").Value = -35.483870967741
$ws.Range("M28").Value = -15.966386554621
$ws.Range("N28").Value = -79.757085020242
$ws.Range("C29").Value = 3
$ws.Range("E29").Value = 50
$ws.Range("F29").Value = 17
$ws.Range("H29").Value = 54.545454545454
$ws.Range("I29").Value = 89
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = -11
$ws.Range("L29").Value = -35.036496350365
$ws.Range("M29").Value = -14.423076923076
$ws.Range("N29").Value = -80.353200883002
$ws.Range("D30").Value = 3
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 66
$ws.Range("K30").Value = -36.363636363636
$ws.Range("L30").Value = -14.285714285714
